$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6496
$ws.Range("C26").Value = 1012
$ws.Range("D26").Value = 6056536
$ws.Range("E26").Value = 932.3485221674877
$ws.Range("F26").Value = 9.785364204833535
$ws.Range("G26").Value = 7.430997876857748
$ws.Range("H26").Value = 26.12756624287029
